$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Collapse the three "CORE COMPETENCIES" bullet paragraphs into a
#    single summary paragraph.
# ------------------------------------------------------------------
# Paragraph 6 = "Product Marketing Core: ..."
# Paragraph 7 = "Research & Analytics: ..."
# Paragraph 8 = "Communication & Technology: ..."
$p6 = $d.Paragraphs.Item(6)

# Delete paragraphs 7 and 8 entirely (including their paragraph marks),
# leaving only paragraph 6 which we then rewrite. Re-fetch paragraph 7
# after each delete since the collection re-indexes live.
$d.Paragraphs.Item(7).Range.Delete()
$d.Paragraphs.Item(7).Range.Delete()

$p6.Range.Text = "Product Marketing Core " + [char]0x2022 + " Research & Analytics " + [char]0x2022 + " Communication & Technology"

# ------------------------------------------------------------------
# 2. Append a new "TECHNICAL SKILLS" section at the end of the
#    document (before the section properties), containing the
#    original detailed skills text that used to live under
#    CORE COMPETENCIES.
# ------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($lastIndex)

# Insert four new paragraphs after the final bullet point while the
# style is still inherited as "Normal" from that last bullet, so the
# three body paragraphs end up with no explicit pPr/pStyle.
$last.Range.InsertParagraphAfter()
$headingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$headingPara.Range.Text = "TECHNICAL SKILLS"

$headingPara.Range.InsertParagraphAfter()
$para1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$para1.Range.Text = "PRODUCT MARKETING CORE Market Intelligence & Competitive Analysis; Product Positioning & Messaging Development; Go-to-Market Strategy & Product Launch Management; Customer Segmentation & Buyer Persona Development; Cross-functional Team Leadership & Collaboration; Sales Enablement & Training Material Development; Data-Driven Decision Making & Analytics Interpretation"

$para1.Range.InsertParagraphAfter()
$para2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$para2.Range.Text = "RESEARCH & ANALYTICS Survey Methodology & Customer Insights; Market Research Design & Implementation; Competitive Intelligence & SWOT Analysis; Customer Journey Mapping & Behavioral Analysis; Statistical Modeling & Trend Analysis; Performance Metrics & Dashboard Development; A/B Testing & Conversion Optimization"

$para2.Range.InsertParagraphAfter()
$para3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$para3.Range.Text = "COMMUNICATION & TECHNOLOGY Strategic Messaging & Narrative Development; Technical Concept Translation for Business Audiences; Stakeholder Communication & Presentation Skills; Data Visualization & Reporting (Tableau, PowerBI, d3.js); Marketing Technology Stack Integration; Content Strategy & Thought Leadership; Client Relationship Management & Business Development"

# Now apply the Heading2 style only to the new section heading, after
# all the body paragraphs have already been created as "Normal".
$headingPara.Style = "Heading2"
